# Apply trade #45 close update across the workbook.
# Trade #45 (row 46 in "All Trades" / "MarketMaking" detail sheets) closed
# at 2026-02-17 15:30:30 with a small loss, which ripples into the
# Summary and Strategy Status roll-up sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.71   # Current Capital
$summary.Range("B4").Value = 0.71      # Total P&L $
$summary.Range("B5").Value = 0.32      # Total P&L %
$summary.Range("B6").Value = 45        # Total Trades
$summary.Range("B8").Value = 22        # Losing Trades
$summary.Range("B9").Value = 31.11     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.71     # Capital
$status.Range("D4").Value = 45         # Trades
$status.Range("E4").Value = 0.71       # P&L $
$status.Range("F4").Value = 0.71       # P&L %
$status.Range("G4").Value = 31.11      # Win Rate %

# ---------------------------------------------------------------
# All Trades + MarketMaking detail sheets - row 46 (Trade #45)
# ---------------------------------------------------------------
$detailSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $detailSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G46").Value = 0.34          # Exit Price
    $ws.Range("H46").Value = "CLOSED"      # Status
    $ws.Range("I46").Value = -10.1297      # P&L %
    $ws.Range("J46").Value = -0.04         # P&L $
    $ws.Range("K46").Value = 100.71        # Capital After
    $ws.Range("P46").Value = "early_exit"  # Exit Reason
    $ws.Range("Q46").Value = 0.13          # Duration (min)
}
